# Updated Excel table with queries
# - Column A (query text) for rows 24-37 is re-populated with the values
#   that, before this edit, lived 3 rows further down (rows 27-40), each
#   gaining the "Good" (green) cell style that already decorated A1:A23.
# - The trailing rows (A38:A40) that no longer have a corresponding query
#   are cleared out entirely.
# - The active selection/zoom of the sheet is updated to match the new
#   view state, and column A's width is nudged back to a round 61
#   characters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column-A text for rows 24-37 (rows 24-28 keep their existing text,
# rows 29-37 take on the text that used to sit 3 rows below them).
$values = @{
    24 = "corona tested positive yesterday europe"
    25 = "Number of positive tests in russia"
    26 = "Last year, how many people tested positive for the coronavirus in Spain alone?"
    27 = "How many vaccinations were performed in Uganda on July 2nd, 2021?"
    28 = "How many vaccinations were performed in Hungary on July 2nd, 2019?"
    29 = "How many new cases were discovered in Serbia in 2021?"
    30 = "How many new cases were detected in Russia in 2018?"
    31 = "What is the peak number of confirmed cases in certain country"
    32 = "What is the peak number of vaccinated people in a day of a certain country"
    33 = "How many new cases of COVID are there today in Hong Kong?"
    34 = "How many new cases of COVID are this week in Hong Kong?"
    35 = "Which country has had the most corona cases?"
    36 = "On which day were the most cases reported?"
    37 = "In which country did most people get vaccinated?"
}

foreach ($r in $values.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $values[$r]
    $cell.Style = "Good"
}

# Rows 38-40 no longer have a matching query in column A.
[void]$ws.Range("A38:A40").ClearContents()

# Column A width: round to 61 characters.
$ws.Columns.Item(1).ColumnWidth = 60.17

# Update the view: selection moved to A33:A37, zoomed to 125%, no frozen
# top-left-cell scroll offset.
[void]$ws.Range("A33:A37").Select()
$excel.ActiveWindow.Zoom = 125
